$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("W2").Value = 2.22
$ws.Range("I3").Value = 6.4
$ws.Range("S3").Value = 3.55
$ws.Range("V3").Value = 1.18
$ws.Range("Q4").Value = 2.16
$ws.Range("AA5").Value = 14
$ws.Range("AI5").Value = 32
$ws.Range("I5").Value = 1.56
$ws.Range("R5").Value = 1.46
$ws.Range("S5").Value = 3
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 2.78
$ws.Range("AC6").Value = 17
$ws.Range("AE6").Value = 170
$ws.Range("AI6").Value = 130
$ws.Range("AK6").Value = 12.5
$ws.Range("F6").Value = 1.26
$ws.Range("G6").Value = 1.27
$ws.Range("I6").Value = 13.5
$ws.Range("K6").Value = 7.8
$ws.Range("R6").Value = 1.95
$ws.Range("V6").Value = 1.08
$ws.Range("W6").Value = 4.7
$ws.Range("AB7").Value = 16.5
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 27
$ws.Range("AG7").Value = 14.5
$ws.Range("AI7").Value = 30
$ws.Range("AK7").Value = 38
$ws.Range("AM7").Value = 70
$ws.Range("AN7").Value = 29
$ws.Range("H7").Value = 2.14
$ws.Range("I7").Value = 2.2
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 3.85
$ws.Range("L7").Value = 1.33
$ws.Range("N7").Value = 4.8
$ws.Range("P7").Value = 2.28
$ws.Range("Q7").Value = 1.74
$ws.Range("S7").Value = 2.84
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 1.37
$ws.Range("X7").Value = 18.5
$ws.Range("Z7").Value = 15
$ws.Range("AA8").Value = 320
$ws.Range("AC8").Value = 10.5
$ws.Range("AE8").Value = 150
$ws.Range("AG8").Value = 9.800000000000001
$ws.Range("AH8").Value = 26
$ws.Range("AJ8").Value = 12.5
$ws.Range("AL8").Value = 38
$ws.Range("AM8").Value = 180
$ws.Range("AN8").Value = 7.2
$ws.Range("AO8").Value = 190
$ws.Range("H8").Value = 7.6
$ws.Range("J8").Value = 5.1
$ws.Range("L8").Value = 1.35
$ws.Range("T8").Value = 2.06
$ws.Range("V8").Value = 1.13
$ws.Range("W8").Value = 3.05
$ws.Range("X8").Value = 18
$ws.Range("AC9").Value = 9.4
$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 19.5
$ws.Range("AF9").Value = 27
$ws.Range("AG9").Value = 14
$ws.Range("AH9").Value = 14.5
$ws.Range("AI9").Value = 27
$ws.Range("AK9").Value = 30
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 55
$ws.Range("AN9").Value = 20
$ws.Range("H9").Value = 2.22
$ws.Range("J9").Value = 4.1
$ws.Range("L9").Value = 1.28
$ws.Range("V9").Value = 1.8
$ws.Range("W9").Value = 1.44
$ws.Range("X9").Value = 24
$ws.Range("F10").Value = 2.28
$ws.Range("I10").Value = 3.2
$ws.Range("L10").Value = 1.28
$ws.Range("S10").Value = 2.4
$ws.Range("U10").Value = 2.74
$ws.Range("V10").Value = 1.45
$ws.Range("W10").Value = 1.76
$ws.Range("Z10").Value = 27
$ws.Range("AC11").Value = 8.6
$ws.Range("AD11").Value = 13.5
$ws.Range("AF11").Value = 16
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 15
$ws.Range("AK11").Value = 21
$ws.Range("AL11").Value = 30
$ws.Range("AO11").Value = 25
$ws.Range("F11").Value = 2.26
$ws.Range("G11").Value = 2.28
$ws.Range("J11").Value = 3.8
$ws.Range("L11").Value = 1.33
$ws.Range("P11").Value = 2.32
$ws.Range("S11").Value = 2.8
$ws.Range("V11").Value = 1.41
$ws.Range("W11").Value = 1.78
$ws.Range("X11").Value = 18.5
$ws.Range("Y11").Value = 16.5
$ws.Range("Z11").Value = 25
$ws.Range("AB12").Value = 19.5
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 240
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 34
$ws.Range("J12").Value = 9.800000000000001
$ws.Range("L12").Value = 1.15
$ws.Range("T12").Value = 1.82
$ws.Range("U12").Value = 2.12
$ws.Range("V12").Value = 1.05
$ws.Range("W12").Value = 6.2
$ws.Range("Y12").Value = 95
$ws.Range("Z12").Value = 230
$ws.Range("AA13").Value = 32
$ws.Range("AD13").Value = 11.5
$ws.Range("AE13").Value = 22
$ws.Range("AF13").Value = 24
$ws.Range("AG13").Value = 13
$ws.Range("AH13").Value = 14.5
$ws.Range("AI13").Value = 29
$ws.Range("AK13").Value = 30
$ws.Range("AL13").Value = 36
$ws.Range("AM13").Value = 60
$ws.Range("AN13").Value = 21
$ws.Range("AO13").Value = 13.5
$ws.Range("J13").Value = 3.75
$ws.Range("L13").Value = 1.29
$ws.Range("P13").Value = 2.38
$ws.Range("V13").Value = 1.71
$ws.Range("W13").Value = 1.46
$ws.Range("X13").Value = 21
$ws.Range("Z13").Value = 17.5
